$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Overview sheet: the file has been handed back, update the status text.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Helper data reused for both the zh-cn and de-de detail sheets.
# ---------------------------------------------------------------------------
$mdFile      = "20835d32-3b2b-4e1f-9dff-9cf50e9033c2.md"
$mdUrl       = "https://github.com/OpenLocalizationTest/oltest/blob/d3b3b3ffaad2419d854271a92ec72658690b58bf/e2e/20835d32-3b2b-4e1f-9dff-9cf50e9033c2.md"
$ffffFile    = "ffff328d9201-3368-4723-993e-08ce926f2040.md"
$ffffUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/d3b3b3ffaad2419d854271a92ec72658690b58bf/e2e/ffff328d9201-3368-4723-993e-08ce926f2040.md"

# ---------------------------------------------------------------------------
# 2) zh-cn sheet
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$zhcnXlf    = "20835d32-3b2b-4e1f-9dff-9cf50e9033c2.ec130cb566f2e5bc09968064a9507e6a1f74308f.zh-cn.xlf"
$zhcnXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5d59be7e7b1ccb66dbf7e25294f639de37b99cf8/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/20835d32-3b2b-4e1f-9dff-9cf50e9033c2.ec130cb566f2e5bc09968064a9507e6a1f74308f.zh-cn.xlf"
$zhcnHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5d59be7e7b1ccb66dbf7e25294f639de37b99cf8/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/hb/20835d32-3b2b-4e1f-9dff-9cf50e9033c2.ec130cb566f2e5bc09968064a9507e6a1f74308f.zh-cn.xlf"

# New "Latest Target File" / "Latest Handback File" columns for rows 2 and 3.
$ws2.Range("F2").Value = $mdFile
$ws2.Range("F3").Value = $mdFile
$ws2.Range("G2").Value = $zhcnXlf
$ws2.Range("G3").Value = $zhcnXlf

# Give the new cells the same look as the other file-name / hyperlink cells
# in the table (underlined custom-blue font, matching style index 1).
$ws2.Range("F2:G3").Font.Underline = 2
$ws2.Range("F2:G3").Font.Color = 15570276
$ws2.Range("F2:G3").Font.Name = "Calibri"
$ws2.Range("F2:G3").Font.Size = 11

# The file has now actually been handed back -> stamp the handback datetime.
$ws2.Range("H2").Value = "2016-03-22 05:43:38"
$ws2.Range("H3").Value = "2016-03-22 05:43:38"

# Rebuild the hyperlinks collection so the new links land in worksheet
# (row-major) order, right after the existing ones for each row.
$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $mdUrl, "", "", $mdFile)
$ws2.Hyperlinks.Add($ws2.Range("D2"), $zhcnXlfUrl, "", "", $zhcnXlf)
$ws2.Hyperlinks.Add($ws2.Range("F2"), $mdUrl, "", "", $mdFile)
$ws2.Hyperlinks.Add($ws2.Range("G2"), $zhcnHandbackUrl, "", "", $zhcnXlf)
$ws2.Hyperlinks.Add($ws2.Range("A3"), $ffffUrl, "", "", $ffffFile)
$ws2.Hyperlinks.Add($ws2.Range("D3"), $zhcnXlfUrl, "", "", $zhcnXlf)
$ws2.Hyperlinks.Add($ws2.Range("F3"), $mdUrl, "", "", $mdFile)
$ws2.Hyperlinks.Add($ws2.Range("G3"), $zhcnHandbackUrl, "", "", $zhcnXlf)

# ---------------------------------------------------------------------------
# 3) de-de sheet
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$dedeXlf    = "20835d32-3b2b-4e1f-9dff-9cf50e9033c2.ec130cb566f2e5bc09968064a9507e6a1f74308f.de-de.xlf"
$dedeXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e3891e992c3c0ec8e759739ddc6ac9b967510215/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/20835d32-3b2b-4e1f-9dff-9cf50e9033c2.ec130cb566f2e5bc09968064a9507e6a1f74308f.de-de.xlf"
$dedeHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e3891e992c3c0ec8e759739ddc6ac9b967510215/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/hb/20835d32-3b2b-4e1f-9dff-9cf50e9033c2.ec130cb566f2e5bc09968064a9507e6a1f74308f.de-de.xlf"

$ws3.Range("F2").Value = $mdFile
$ws3.Range("F3").Value = $mdFile
$ws3.Range("G2").Value = $dedeXlf
$ws3.Range("G3").Value = $dedeXlf

$ws3.Range("F2:G3").Font.Underline = 2
$ws3.Range("F2:G3").Font.Color = 15570276
$ws3.Range("F2:G3").Font.Name = "Calibri"
$ws3.Range("F2:G3").Font.Size = 11

# A second, independent file, already handed back a bit later.
$ws3.Range("H2").Value = "2016-03-22 05:43:52"
$ws3.Range("H3").Value = "2016-03-22 05:43:52"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $mdUrl, "", "", $mdFile)
$ws3.Hyperlinks.Add($ws3.Range("D2"), $dedeXlfUrl, "", "", $dedeXlf)
$ws3.Hyperlinks.Add($ws3.Range("F2"), $mdUrl, "", "", $mdFile)
$ws3.Hyperlinks.Add($ws3.Range("G2"), $dedeHandbackUrl, "", "", $dedeXlf)
$ws3.Hyperlinks.Add($ws3.Range("A3"), $ffffUrl, "", "", $ffffFile)
$ws3.Hyperlinks.Add($ws3.Range("D3"), $dedeXlfUrl, "", "", $dedeXlf)
$ws3.Hyperlinks.Add($ws3.Range("F3"), $mdUrl, "", "", $mdFile)
$ws3.Hyperlinks.Add($ws3.Range("G3"), $dedeHandbackUrl, "", "", $dedeXlf)
